$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.400.11"
$ws.Range("E2").Value = "  +0.12%  "
$ws.Range("D3").Value = "1.849.47"
$ws.Range("E3").Value = "  +0.19%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'240.58"
$ws.Range("E5").Value = "  +0.21%  "
$ws.Range("D6").Value = "'0.6301"
$ws.Range("E6").Value = "  -0.10%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "'0.07659"
$ws.Range("E8").Value = "  +1.61%  "
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("D10").Value = "'24.48"
$ws.Range("E10").Value = "  +0.05%  "
$ws.Range("D11").Value = "'0.07742"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").Value = "1.842.43"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").Value = "'5.011"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").Value = "'0.00001092"
$ws.Range("E14").Value = "  +9.20%  "
$ws.Range("D15").Value = "'0.6796"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16").Value = "'83.58"
$ws.Range("E16").Value = "  +0.90%  "
$ws.Range("D17").Value = "2.088.74"
$ws.Range("E17").Value = "  -7.76%  "
$ws.Range("D18").Value = "'6.139"
$ws.Range("E18").Value = "  -0.04%  "
$ws.Range("D19").Value = "29.421.77"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "'229.59"
$ws.Range("E20").Value = "  +0.82%  "
$ws.Range("D21").Value = "'12.45"
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'7.439"
$ws.Range("E23").Value = "  -1.34%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("D25").Value = "'157.15"
$ws.Range("D26").Value = "'0.1390"
$ws.Range("E26").Value = "  -0.55%  "
$ws.Range("D27").Value = "'8.372"
$ws.Range("E27").Value = "  +0.02%  "
$ws.Range("E28").Value = "  +0.15%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("D30").Value = "'1.305"
$ws.Range("E30").Value = "  +3.94%  "
$ws.Range("D31").Value = "'0.05676"
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("E32").Value = "  -0.33%  "
$ws.Range("D33").Value = "'4.045"
$ws.Range("E33").Value = "  +0.80%  "
$ws.Range("E34").Value = "  +0.36%  "
$ws.Range("E35").Value = "  +0.24%  "
$ws.Range("D36").Value = "'0.7097"
$ws.Range("E36").Value = "  -0.62%  "
$ws.Range("D37").Value = "'2.585"
$ws.Range("D38").Value = "'2.776"
$ws.Range("E38").Value = "  -0.26%  "
$ws.Range("D39").Value = "1.229.60"
$ws.Range("E39").Value = "  -2.11%  "
$ws.Range("D40").Value = "'0.01798"
$ws.Range("E40").Value = "  -0.80%  "
$ws.Range("D41").Value = "'6.467"
$ws.Range("E41").Value = "  +4.24%  "
$ws.Range("D42").Value = "'0.9140"
$ws.Range("E42").Value = "  -0.02%  "
$ws.Range("E43").Value = "  -0.01%  "
$ws.Range("D44").Value = "1.998.23"
$ws.Range("E44").Value = "  -0.18%  "
$ws.Range("D45").Value = "'101.42"
$ws.Range("E45").Value = "  +0.15%  "
$ws.Range("D46").Value = "'66.18"
$ws.Range("E46").Value = "  -0.33%  "
$ws.Range("E47").Value = "  +4.03%  "
$ws.Range("D48").Value = "'7.154"
$ws.Range("E48").Value = "  +1.52%  "
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").Value = "'9.036"
$ws.Range("E50").Value = "  -0.89%  "
$ws.Range("E51").Value = "  -0.22%  "
